$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the schema validation flag from true to false in cell J3
$ws.Range("J3").Value = "pets_json.schema.json=false"

# Update the selected/active cell on the sheet to J18
$ws.Range("J18").Select()
